$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 240.33333
$ws.Range("I4").Value = 133.33333
$ws.Range("J4").Value = 347.33334
$ws.Range("K4").Value = 133.33333
$ws.Range("L4").Value = 347.33334
$ws.Range("M4").Value = -19.33332999999999
$ws.Range("N4").Value = -575.33334
$ws.Range("H48").Value = 12032.6
$ws.Range("J48").Value = 13721
$ws.Range("L48").Value = 41163
$ws.Range("N48").Value = -41747
$ws.Range("H56").Value = 12032.6
$ws.Range("J56").Value = 13721
$ws.Range("L56").Value = 41163
$ws.Range("N56").Value = -42231
$ws.Range("H76").Value = 21935.25
$ws.Range("I76").Value = 29397.4
$ws.Range("K76").Value = 29397.4
$ws.Range("M76").Value = -29082.4
$ws.Range("H79").Value = 21935.25
$ws.Range("I79").Value = 29397.4
$ws.Range("K79").Value = 29397.4
$ws.Range("M79").Value = -28305.4
$ws.Range("H112").Value = 1392469.9
$ws.Range("J112").Value = 1591056
$ws.Range("L112").Value = 4773168
$ws.Range("N112").Value = -4775384
$ws.Range("H113").Value = 4994.1934
$ws.Range("I113").Value = 6980.1055
$ws.Range("K113").Value = 6980.1055
$ws.Range("M113").Value = -3726.1055
$ws.Range("H132").Value = 1663.193
$ws.Range("I132").Value = 1666.12
$ws.Range("K132").Value = 4998.36
$ws.Range("M132").Value = -2468.36
$ws.Range("H137").Value = 1794207
$ws.Range("I137").Value = 1909.6471
$ws.Range("J137").Value = 3970568
$ws.Range("K137").Value = 5728.9413
$ws.Range("L137").Value = 11911704
$ws.Range("M137").Value = -3178.9413
$ws.Range("N137").Value = -11916804
$ws.Range("H138").Value = 2308.192
$ws.Range("J138").Value = 2739.9734
$ws.Range("L138").Value = 8219.9202
$ws.Range("N138").Value = -18499.9202

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3216.7334
$ws.Range("I2").Value = 2861.4443
$ws.Range("K2").Value = 2861.4443
$ws.Range("M2").Value = -2748.4443
$ws.Range("H45").Value = 2930.8333
$ws.Range("I45").Value = 1217.5
$ws.Range("K45").Value = 1217.5
$ws.Range("M45").Value = -840.5
$ws.Range("H102").Value = 2424.5334
$ws.Range("I102").Value = 1997.1818
$ws.Range("K102").Value = 1997.1818
$ws.Range("M102").Value = -375.1818000000001
$ws.Range("H116").Value = 3216.7334
$ws.Range("I116").Value = 2861.4443
$ws.Range("K116").Value = 2861.4443
$ws.Range("M116").Value = -567.4443000000001
$ws.Range("H122").Value = 3671.2222
$ws.Range("I122").Value = 3269.3
$ws.Range("K122").Value = 9807.900000000001
$ws.Range("M122").Value = -7357.900000000001
$ws.Range("H132").Value = 3949.375
$ws.Range("I132").Value = 3050.7585
$ws.Range("K132").Value = 9152.2755
$ws.Range("M132").Value = -6622.2755

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3216.7334
$ws.Range("I3").Value = 2861.4443
$ws.Range("K3").Value = 2861.4443
$ws.Range("M3").Value = -2747.4443
$ws.Range("H20").Value = 43026.668
$ws.Range("I20").Value = 62851.375
$ws.Range("K20").Value = 62851.375
$ws.Range("M20").Value = -62604.375
$ws.Range("H134").Value = 1834755.9
$ws.Range("I134").Value = 2383009.5
$ws.Range("J134").Value = 7243.6665
$ws.Range("K134").Value = 7149028.5
$ws.Range("L134").Value = 21730.9995
$ws.Range("M134").Value = -7146493.5
$ws.Range("N134").Value = -26800.9995

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5443.107
$ws.Range("I31").Value = 2200.4119
$ws.Range("K31").Value = 2200.4119
$ws.Range("M31").Value = -1905.4119
$ws.Range("H34").Value = 5443.107
$ws.Range("I34").Value = 2200.4119
$ws.Range("K34").Value = 2200.4119
$ws.Range("M34").Value = -1998.4119
$ws.Range("H35").Value = 676666.3
$ws.Range("I35").Value = 676666.3
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 676666.3
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -676372.3
$ws.Range("N35").ClearContents()
$ws.Range("H64").Value = 67663.336
$ws.Range("J64").Value = 67663.336
$ws.Range("L64").Value = 67663.336
$ws.Range("N64").Value = -68159.336
$ws.Range("H67").Value = 67663.336
$ws.Range("J67").Value = 67663.336
$ws.Range("L67").Value = 67663.336
$ws.Range("N67").Value = -69379.336
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 1258.36
$ws.Range("I107").Value = 593.4167
$ws.Range("K107").Value = 593.4167
$ws.Range("M107").Value = 1326.5833
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4016
$ws.Range("I132").Value = 4032.389
$ws.Range("J132").Value = 3957
$ws.Range("K132").Value = 12097.167
$ws.Range("L132").Value = 11871
$ws.Range("M132").Value = -9567.167000000001
$ws.Range("N132").Value = -16931
$ws.Range("H134").Value = 3040.1
$ws.Range("I134").Value = 3040.1
$ws.Range("K134").Value = 9120.299999999999
$ws.Range("M134").Value = -6585.299999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 197.5
$ws.Range("I87").Value = 197.5
$ws.Range("K87").Value = 592.5
$ws.Range("M87").Value = 655.5
$ws.Range("H90").Value = 197.5
$ws.Range("I90").Value = 197.5
$ws.Range("K90").Value = 1777.5
$ws.Range("M90").Value = 4462.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 38400.6
$ws.Range("I20").Value = 12005
$ws.Range("J20").Value = 44999.5
$ws.Range("K20").Value = 12005
$ws.Range("L20").Value = 44999.5
$ws.Range("M20").Value = -11760
$ws.Range("N20").Value = -45489.5
$ws.Range("H70").Value = 18452.355
$ws.Range("I70").Value = 55904
$ws.Range("J70").Value = 4833.5757
$ws.Range("K70").Value = 55904
$ws.Range("L70").Value = 4833.5757
$ws.Range("M70").Value = -55634
$ws.Range("N70").Value = -5373.5757
$ws.Range("H73").Value = 18452.355
$ws.Range("I73").Value = 55904
$ws.Range("J73").Value = 4833.5757
$ws.Range("K73").Value = 55904
$ws.Range("L73").Value = 4833.5757
$ws.Range("M73").Value = -54968
$ws.Range("N73").Value = -6705.5757
$ws.Range("H102").Value = 3068.5
$ws.Range("I102").Value = 2910.625
$ws.Range("J102").Value = 3700
$ws.Range("K102").Value = 2910.625
$ws.Range("L102").Value = 3700
$ws.Range("M102").Value = -1288.625
$ws.Range("N102").Value = -6944
$ws.Range("H122").Value = 1996.8334
$ws.Range("I122").Value = 1624.5
$ws.Range("K122").Value = 4873.5
$ws.Range("M122").Value = -2423.5
$ws.Range("H132").Value = 4268.769
$ws.Range("I132").Value = 3549.5
$ws.Range("K132").Value = 10648.5
$ws.Range("M132").Value = -8118.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 319.92856
$ws.Range("I16").Value = 329.15384
$ws.Range("K16").Value = 329.15384
$ws.Range("M16").Value = -159.15384
$ws.Range("H82").Value = 3851
$ws.Range("I82").Value = 3978.5
$ws.Range("K82").Value = 3978.5
$ws.Range("M82").Value = -3617.5
$ws.Range("H85").Value = 3851
$ws.Range("I85").Value = 3978.5
$ws.Range("K85").Value = 3978.5
$ws.Range("M85").Value = -2730.5
$ws.Range("H132").Value = 3283.6
$ws.Range("J132").Value = 3999.5
$ws.Range("L132").Value = 11998.5
$ws.Range("N132").Value = -17058.5
$ws.Range("H136").Value = 2707.8462
$ws.Range("I136").Value = 2437.25
$ws.Range("J136").Value = 3140.8
$ws.Range("K136").Value = 7311.75
$ws.Range("L136").Value = 9422.400000000001
$ws.Range("M136").Value = -4761.75
$ws.Range("N136").Value = -14522.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 30015.25
$ws.Range("J21").Value = 30015.25
$ws.Range("L21").Value = 30015.25
$ws.Range("N21").Value = -30485.25
$ws.Range("H35").Value = 30015.25
$ws.Range("J35").Value = 30015.25
$ws.Range("L35").Value = 30015.25
$ws.Range("N35").Value = -30595.25
$ws.Range("H39").Value = 28000
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H132").Value = 3435.7368
$ws.Range("I132").Value = 3411.25
$ws.Range("J132").Value = 3566.3333
$ws.Range("K132").Value = 10233.75
$ws.Range("L132").Value = 10698.9999
$ws.Range("M132").Value = -7703.75
$ws.Range("N132").Value = -15758.9999
